$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values updated with the latest automatic electricity price feed.
$ws.Range("A2").Value = 45891
$ws.Range("B2").Value = 77.09999999999999
$ws.Range("C2").Value = 72
$ws.Range("D2").Value = 68.01000000000001
$ws.Range("E2").Value = 65.2
$ws.Range("F2").Value = 70.01000000000001
$ws.Range("G2").Value = 69.81
$ws.Range("H2").Value = 71.03
$ws.Range("I2").Value = 78.06
$ws.Range("J2").Value = 75.09999999999999
$ws.Range("K2").Value = 49.67
$ws.Range("L2").Value = 15.57
$ws.Range("M2").Value = 5.11
$ws.Range("N2").Value = 4.31
$ws.Range("O2").Value = 5.79
$ws.Range("P2").Value = 5.37
$ws.Range("Q2").Value = 5.01
$ws.Range("R2").Value = 5.79
$ws.Range("S2").Value = 20
$ws.Range("T2").Value = 35.2
$ws.Range("U2").Value = 84.56
$ws.Range("V2").Value = 104.53
$ws.Range("W2").Value = 103.49
$ws.Range("X2").Value = 101.59
$ws.Range("Y2").Value = 94.37
$ws.Range("Z2").Value = 53.61
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 101
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 104.01
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 97.98
$ws.Range("AG2").Value = "9h-18h"
